$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=2; Col=4; Value='69.036.74'},
    @{Row=2; Col=5; Value='  +2.04%  '},
    @{Row=3; Col=4; Value='2.523.86'},
    @{Row=3; Col=5; Value='  +1.45%  '},
    @{Row=4; Col=4; Value='1.00'},
    @{Row=4; Col=5; Value='  +0.10%  '},
    @{Row=5; Col=4; Value='595.23'},
    @{Row=5; Col=5; Value='  +1.61%  '},
    @{Row=6; Col=4; Value='176.23'},
    @{Row=6; Col=5; Value='  -0.06%  '},
    @{Row=7; Col=5; Value='  +0.04%  '},
    @{Row=8; Col=4; Value='0.519'},
    @{Row=8; Col=5; Value='  +0.82%  '},
    @{Row=9; Col=4; Value='2.523.33'},
    @{Row=9; Col=5; Value='  +1.46%  '},
    @{Row=10; Col=4; Value='0.148'},
    @{Row=10; Col=5; Value='  +6.23%  '},
    @{Row=11; Col=5; Value='  -1.05%  '},
    @{Row=12; Col=5; Value='  +1.25%  '},
    @{Row=13; Col=4; Value='0.339'},
    @{Row=13; Col=5; Value='  +0.50%  '},
    @{Row=14; Col=4; Value='26.17'},
    @{Row=14; Col=5; Value='  +1.70%  '},
    @{Row=15; Col=4; Value='2.944.99'},
    @{Row=15; Col=5; Value='  -0.17%  '},
    @{Row=16; Col=4; Value='68.640.69'},
    @{Row=16; Col=5; Value='  +1.58%  '},
    @{Row=17; Col=5; Value='  +1.09%  '},
    @{Row=18; Col=4; Value='2.528.79'},
    @{Row=18; Col=5; Value='  +2.91%  '},
    @{Row=19; Col=4; Value='11.09'},
    @{Row=19; Col=5; Value='  +0.33%  '},
    @{Row=20; Col=4; Value='362.38'},
    @{Row=20; Col=5; Value='  +3.15%  '},
    @{Row=21; Col=4; Value='7.52'},
    @{Row=21; Col=5; Value='  +1.19%  '},
    @{Row=22; Col=4; Value='4.10'},
    @{Row=22; Col=5; Value='  +0.61%  '},
    @{Row=23; Col=5; Value='  +0.06%  '},
    @{Row=24; Col=4; Value='70.65'},
    @{Row=24; Col=5; Value='  +0.02%  '},
    @{Row=25; Col=4; Value='4.21'},
    @{Row=25; Col=5; Value='  -0.26%  '},
    @{Row=26; Col=4; Value='1.70'},
    @{Row=26; Col=5; Value='  -5.21%  '},
    @{Row=27; Col=4; Value='9.04'},
    @{Row=27; Col=5; Value='  -2.26%  '},
    @{Row=28; Col=5; Value='  +1.52%  '},
    @{Row=29; Col=4; Value='0.995'},
    @{Row=29; Col=5; Value='  -0.33%  '},
    @{Row=30; Col=4; Value='520.60'},
    @{Row=30; Col=5; Value='  +2.76%  '},
    @{Row=31; Col=4; Value='0.0₃0894'},
    @{Row=31; Col=5; Value='  -1.72%  '},
    @{Row=32; Col=4; Value='7.78'},
    @{Row=32; Col=5; Value='  -0.64%  '},
    @{Row=33; Col=5; Value='  +0.17%  '},
    @{Row=34; Col=5; Value='  +0.27%  '},
    @{Row=35; Col=4; Value='1.00'},
    @{Row=35; Col=5; Value='  +0.05%  '},
    @{Row=36; Col=4; Value='163.20'},
    @{Row=36; Col=5; Value='  +1.39%  '},
    @{Row=37; Col=5; Value='  -1.57%  '},
    @{Row=38; Col=4; Value='18.54'},
    @{Row=38; Col=5; Value='  +1.19%  '},
    @{Row=39; Col=4; Value='18.69'},
    @{Row=40; Col=4; Value='1.77'},
    @{Row=40; Col=5; Value='  +3.50%  '},
    @{Row=41; Col=5; Value='  -1.09%  '},
    @{Row=42; Col=5; Value='  +0.04%  '},
    @{Row=43; Col=4; Value='4.84'},
    @{Row=43; Col=5; Value='  -0.43%  '},
    @{Row=44; Col=4; Value='0.326'},
    @{Row=44; Col=5; Value='  -0.94%  '},
    @{Row=45; Col=4; Value='2.40'},
    @{Row=45; Col=5; Value='  -1.25%  '},
    @{Row=46; Col=4; Value='151.29'},
    @{Row=46; Col=5; Value='  +5.45%  '},
    @{Row=47; Col=4; Value='3.59'},
    @{Row=47; Col=5; Value='  +2.40%  '},
    @{Row=48; Col=4; Value='0.518'},
    @{Row=48; Col=5; Value='  +1.08%  '},
    @{Row=49; Col=2; Value='Cronos'},
    @{Row=49; Col=3; Value='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'},
    @{Row=49; Col=4; Value='0.0742'},
    @{Row=49; Col=5; Value='  -0.46%  '},
    @{Row=50; Col=2; Value='Optimism'},
    @{Row=50; Col=3; Value='https://coinranking.com/coin/n1p-s_gm1+optimism-op'},
    @{Row=50; Col=4; Value='1.59'},
    @{Row=50; Col=5; Value='  +0.75%  '},
    @{Row=51; Col=4; Value='0.580'},
    @{Row=51; Col=5; Value='  -1.03%  '}
)

foreach ($chg in $changes) {
    $cell = $ws.Cells.Item($chg.Row, $chg.Col)
    $cell.Value = "'" + $chg.Value
    $cell.Style = "Normal"
}
